$wb = $excel.ActiveWorkbook

# ---- Sheet: ALERTS ----
$ws = $wb.Worksheets.Item('ALERTS')
$rng = $ws.Range("A4:F4")
$rng.NumberFormat = "@"
$ws.Range("A4").Value = '2026-01-30'
$ws.Range("B4").Value = '18:19:42'
$ws.Range("C4").Value = '18:00'
$ws.Range("D4").Value = 'Living Room'
$ws.Range("E4").Value = 'CRITICAL'
$ws.Range("F4").Value = 'FALL_DETECTED'

$rng = $ws.Range("A5:F5")
$rng.NumberFormat = "@"
$ws.Range("A5").Value = '2026-01-30'
$ws.Range("B5").Value = '18:19:45'
$ws.Range("C5").Value = '18:00'
$ws.Range("D5").Value = 'Living Room'
$ws.Range("E5").Value = 'CRITICAL'
$ws.Range("F5").Value = 'FALL_DETECTED'


# ---- Sheet: PIR ----
$ws = $wb.Worksheets.Item('PIR')
$rng = $ws.Range("A14:F14")
$rng.NumberFormat = "@"
$ws.Range("A14").Value = '2026-01-30'
$ws.Range("B14").Value = '18:17:15'
$ws.Range("C14").Value = '18:00'
$ws.Range("D14").Value = 'Bathroom'
$ws.Range("E14").Value = 'No Motion'
$ws.Range("F14").Value = 'Inactive'

$rng = $ws.Range("A15:F15")
$rng.NumberFormat = "@"
$ws.Range("A15").Value = '2026-01-30'
$ws.Range("B15").Value = '18:17:15'
$ws.Range("C15").Value = '18:00'
$ws.Range("D15").Value = 'Bathroom'
$ws.Range("E15").Value = 'No Motion'
$ws.Range("F15").Value = 'Inactive'

$rng = $ws.Range("A16:F16")
$rng.NumberFormat = "@"
$ws.Range("A16").Value = '2026-01-30'
$ws.Range("B16").Value = '18:17:20'
$ws.Range("C16").Value = '18:00'
$ws.Range("D16").Value = 'Bathroom'
$ws.Range("E16").Value = 'No Motion'
$ws.Range("F16").Value = 'Inactive'

$rng = $ws.Range("A17:F17")
$rng.NumberFormat = "@"
$ws.Range("A17").Value = '2026-01-30'
$ws.Range("B17").Value = '18:19:45'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Bathroom'
$ws.Range("E17").Value = 'No Motion'
$ws.Range("F17").Value = 'Inactive'

$rng = $ws.Range("A18:F18")
$rng.NumberFormat = "@"
$ws.Range("A18").Value = '2026-01-30'
$ws.Range("B18").Value = '18:19:46'
$ws.Range("C18").Value = '18:00'
$ws.Range("D18").Value = 'Bathroom'
$ws.Range("E18").Value = 'No Motion'
$ws.Range("F18").Value = 'Inactive'

$rng = $ws.Range("A19:F19")
$rng.NumberFormat = "@"
$ws.Range("A19").Value = '2026-01-30'
$ws.Range("B19").Value = '18:19:51'
$ws.Range("C19").Value = '18:00'
$ws.Range("D19").Value = 'Bathroom'
$ws.Range("E19").Value = 'No Motion'
$ws.Range("F19").Value = 'Inactive'

$rng = $ws.Range("A20:F20")
$rng.NumberFormat = "@"
$ws.Range("A20").Value = '2026-01-30'
$ws.Range("B20").Value = '18:19:56'
$ws.Range("C20").Value = '18:00'
$ws.Range("D20").Value = 'Bathroom'
$ws.Range("E20").Value = 'No Motion'
$ws.Range("F20").Value = 'Inactive'

$rng = $ws.Range("A21:F21")
$rng.NumberFormat = "@"
$ws.Range("A21").Value = '2026-01-30'
$ws.Range("B21").Value = '18:20:01'
$ws.Range("C21").Value = '18:00'
$ws.Range("D21").Value = 'Bathroom'
$ws.Range("E21").Value = 'No Motion'
$ws.Range("F21").Value = 'Inactive'

$rng = $ws.Range("A22:F22")
$rng.NumberFormat = "@"
$ws.Range("A22").Value = '2026-01-30'
$ws.Range("B22").Value = '18:20:06'
$ws.Range("C22").Value = '18:00'
$ws.Range("D22").Value = 'Bathroom'
$ws.Range("E22").Value = 'No Motion'
$ws.Range("F22").Value = 'Inactive'

$rng = $ws.Range("A23:F23")
$rng.NumberFormat = "@"
$ws.Range("A23").Value = '2026-01-30'
$ws.Range("B23").Value = '18:20:11'
$ws.Range("C23").Value = '18:00'
$ws.Range("D23").Value = 'Bathroom'
$ws.Range("E23").Value = 'No Motion'
$ws.Range("F23").Value = 'Inactive'

$rng = $ws.Range("A24:F24")
$rng.NumberFormat = "@"
$ws.Range("A24").Value = '2026-01-30'
$ws.Range("B24").Value = '18:20:16'
$ws.Range("C24").Value = '18:00'
$ws.Range("D24").Value = 'Bathroom'
$ws.Range("E24").Value = 'No Motion'
$ws.Range("F24").Value = 'Inactive'

$rng = $ws.Range("A25:F25")
$rng.NumberFormat = "@"
$ws.Range("A25").Value = '2026-01-30'
$ws.Range("B25").Value = '18:20:21'
$ws.Range("C25").Value = '18:00'
$ws.Range("D25").Value = 'Bathroom'
$ws.Range("E25").Value = 'No Motion'
$ws.Range("F25").Value = 'Inactive'

$rng = $ws.Range("A26:F26")
$rng.NumberFormat = "@"
$ws.Range("A26").Value = '2026-01-30'
$ws.Range("B26").Value = '18:20:26'
$ws.Range("C26").Value = '18:00'
$ws.Range("D26").Value = 'Bathroom'
$ws.Range("E26").Value = 'No Motion'
$ws.Range("F26").Value = 'Inactive'

$rng = $ws.Range("A27:F27")
$rng.NumberFormat = "@"
$ws.Range("A27").Value = '2026-01-30'
$ws.Range("B27").Value = '18:20:31'
$ws.Range("C27").Value = '18:00'
$ws.Range("D27").Value = 'Bathroom'
$ws.Range("E27").Value = 'No Motion'
$ws.Range("F27").Value = 'Inactive'

$rng = $ws.Range("A28:F28")
$rng.NumberFormat = "@"
$ws.Range("A28").Value = '2026-01-30'
$ws.Range("B28").Value = '18:20:36'
$ws.Range("C28").Value = '18:00'
$ws.Range("D28").Value = 'Bathroom'
$ws.Range("E28").Value = 'No Motion'
$ws.Range("F28").Value = 'Inactive'


# ---- Sheet: Humidity ----
$ws = $wb.Worksheets.Item('Humidity')
$rng = $ws.Range("A12:F12")
$rng.NumberFormat = "@"
$ws.Range("A12").Value = '2026-01-30'
$ws.Range("B12").Value = '18:17:15'
$ws.Range("C12").Value = '18:00'
$ws.Range("D12").Value = 'Bathroom'
$ws.Range("E12").Value = '86.4%'
$ws.Range("F12").Value = 'Active'

$rng = $ws.Range("A13:F13")
$rng.NumberFormat = "@"
$ws.Range("A13").Value = '2026-01-30'
$ws.Range("B13").Value = '18:17:20'
$ws.Range("C13").Value = '18:00'
$ws.Range("D13").Value = 'Bathroom'
$ws.Range("E13").Value = '86.4%'
$ws.Range("F13").Value = 'Active'

$rng = $ws.Range("A14:F14")
$rng.NumberFormat = "@"
$ws.Range("A14").Value = '2026-01-30'
$ws.Range("B14").Value = '18:19:45'
$ws.Range("C14").Value = '18:00'
$ws.Range("D14").Value = 'Bathroom'
$ws.Range("E14").Value = '86.6%'
$ws.Range("F14").Value = 'Active'

$rng = $ws.Range("A15:F15")
$rng.NumberFormat = "@"
$ws.Range("A15").Value = '2026-01-30'
$ws.Range("B15").Value = '18:19:46'
$ws.Range("C15").Value = '18:00'
$ws.Range("D15").Value = 'Bathroom'
$ws.Range("E15").Value = '86.6%'
$ws.Range("F15").Value = 'Active'

$rng = $ws.Range("A16:F16")
$rng.NumberFormat = "@"
$ws.Range("A16").Value = '2026-01-30'
$ws.Range("B16").Value = '18:19:51'
$ws.Range("C16").Value = '18:00'
$ws.Range("D16").Value = 'Bathroom'
$ws.Range("E16").Value = '86.6%'
$ws.Range("F16").Value = 'Active'

$rng = $ws.Range("A17:F17")
$rng.NumberFormat = "@"
$ws.Range("A17").Value = '2026-01-30'
$ws.Range("B17").Value = '18:20:01'
$ws.Range("C17").Value = '18:00'
$ws.Range("D17").Value = 'Bathroom'
$ws.Range("E17").Value = '86.6%'
$ws.Range("F17").Value = 'Active'

$rng = $ws.Range("A18:F18")
$rng.NumberFormat = "@"
$ws.Range("A18").Value = '2026-01-30'
$ws.Range("B18").Value = '18:20:06'
$ws.Range("C18").Value = '18:00'
$ws.Range("D18").Value = 'Bathroom'
$ws.Range("E18").Value = '86.6%'
$ws.Range("F18").Value = 'Active'

$rng = $ws.Range("A19:F19")
$rng.NumberFormat = "@"
$ws.Range("A19").Value = '2026-01-30'
$ws.Range("B19").Value = '18:20:11'
$ws.Range("C19").Value = '18:00'
$ws.Range("D19").Value = 'Bathroom'
$ws.Range("E19").Value = '86.6%'
$ws.Range("F19").Value = 'Active'

$rng = $ws.Range("A20:F20")
$rng.NumberFormat = "@"
$ws.Range("A20").Value = '2026-01-30'
$ws.Range("B20").Value = '18:20:21'
$ws.Range("C20").Value = '18:00'
$ws.Range("D20").Value = 'Bathroom'
$ws.Range("E20").Value = '86.7%'
$ws.Range("F20").Value = 'Active'

$rng = $ws.Range("A21:F21")
$rng.NumberFormat = "@"
$ws.Range("A21").Value = '2026-01-30'
$ws.Range("B21").Value = '18:20:26'
$ws.Range("C21").Value = '18:00'
$ws.Range("D21").Value = 'Bathroom'
$ws.Range("E21").Value = '85.7%'
$ws.Range("F21").Value = 'Active'

$rng = $ws.Range("A22:F22")
$rng.NumberFormat = "@"
$ws.Range("A22").Value = '2026-01-30'
$ws.Range("B22").Value = '18:20:31'
$ws.Range("C22").Value = '18:00'
$ws.Range("D22").Value = 'Bathroom'
$ws.Range("E22").Value = '86.6%'
$ws.Range("F22").Value = 'Active'


# ---- Sheet: Proximity ----
$ws = $wb.Worksheets.Item('Proximity')
$rng = $ws.Range("A3:F3")
$rng.NumberFormat = "@"
$ws.Range("A3").Value = '2026-01-30'
$ws.Range("B3").Value = '18:19:59'
$ws.Range("C3").Value = '18:00'
$ws.Range("D3").Value = 'Living Room Main Door'
$ws.Range("E3").Value = 'ENTER'
$ws.Range("F3").Value = 'User ENTERED Living Room Main Door'


# ---- Sheet: Camera ----
$ws = $wb.Worksheets.Item('Camera')
$rng = $ws.Range("A3:F3")
$rng.NumberFormat = "@"
$ws.Range("A3").Value = '2026-01-30'
$ws.Range("B3").Value = '18:19:59'
$ws.Range("C3").Value = '18:00'
$ws.Range("D3").Value = 'Living Room Main Door'
$ws.Range("E3").Value = 'Image Captured (ENTER)'
$ws.Range("F3").Value = 'Active'

